# Fix the urls due to the update of wikipedia.
# The crawl that populated this sheet was re-run; four new pages were
# discovered (linked from "游戏") and inserted right before the existing
# "益智游戏" row, pushing the rows below it down by four. The occurrence
# counter for the root "人" row also grew as a result of the re-crawl.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert four new rows before row 160 ("益智游戏") ---
$ws.Rows("160:163").Insert()

# Re-apply the same formatting used by column A of the surrounding rows
# (bold, centered, thin-bordered) to the newly inserted cells.
$ws.Cells.Item(159, 1).Copy($ws.Range("A160:A163"))

# --- Row 160: 多人在线战斗竞技场游戏 ---
$ws.Cells.Item(160, 1).Value = 158
$ws.Cells.Item(160, 2).Value = "游戏"
$ws.Cells.Item(160, 3).Value = "https://zh.wikipedia.org/wiki/%E6%B8%B8%E6%88%8F"
$ws.Cells.Item(160, 4).Value = 159
$ws.Cells.Item(160, 5).Value = "https://zh.wikipedia.org/wiki/%E5%A4%9A%E4%BA%BA%E5%9C%A8%E7%BA%BF%E6%88%98%E6%96%97%E7%AB%9E%E6%8A%80%E5%9C%BA%E6%B8%B8%E6%88%8F"
$ws.Cells.Item(160, 6).Value = "多人在线战斗竞技场游戏"
$ws.Cells.Item(160, 7).Value = 3
$ws.Cells.Item(160, 8).Value = "游戏_游戏"
$ws.Cells.Item(160, 9).Value = 1

# --- Row 161: 魔兽争霸III：混乱之治 ---
$ws.Cells.Item(161, 1).Value = 159
$ws.Cells.Item(161, 2).Value = "游戏"
$ws.Cells.Item(161, 3).Value = "https://zh.wikipedia.org/wiki/%E6%B8%B8%E6%88%8F"
$ws.Cells.Item(161, 4).Value = 160
$ws.Cells.Item(161, 5).Value = "https://zh.wikipedia.org/wiki/%E9%AD%94%E5%85%BD%E4%BA%89%E9%9C%B8III%EF%BC%9A%E6%B7%B7%E4%B9%B1%E4%B9%8B%E6%B2%BB"
$ws.Cells.Item(161, 6).Value = "魔兽争霸III：混乱之治"
$ws.Cells.Item(161, 7).Value = 1
$ws.Cells.Item(161, 8).Value = "游戏_游戏"
$ws.Cells.Item(161, 9).Value = 1

# --- Row 162: 英雄联盟 ---
$ws.Cells.Item(162, 1).Value = 160
$ws.Cells.Item(162, 2).Value = "游戏"
$ws.Cells.Item(162, 3).Value = "https://zh.wikipedia.org/wiki/%E6%B8%B8%E6%88%8F"
$ws.Cells.Item(162, 4).Value = 161
$ws.Cells.Item(162, 5).Value = "https://zh.wikipedia.org/wiki/%E8%8B%B1%E9%9B%84%E8%81%94%E7%9B%9F"
$ws.Cells.Item(162, 6).Value = "英雄联盟"
$ws.Cells.Item(162, 7).Value = 1
$ws.Cells.Item(162, 8).Value = "游戏_游戏"
$ws.Cells.Item(162, 9).Value = 1

# --- Row 163: DotA ---
$ws.Cells.Item(163, 1).Value = 161
$ws.Cells.Item(163, 2).Value = "游戏"
$ws.Cells.Item(163, 3).Value = "https://zh.wikipedia.org/wiki/%E6%B8%B8%E6%88%8F"
$ws.Cells.Item(163, 4).Value = 162
$ws.Cells.Item(163, 5).Value = "https://zh.wikipedia.org/wiki/DotA"
$ws.Cells.Item(163, 6).Value = "DotA"
$ws.Cells.Item(163, 7).Value = 1
$ws.Cells.Item(163, 8).Value = "游戏_游戏"
$ws.Cells.Item(163, 9).Value = 1

# The sequence columns A (0-based) and D (1-based) hold plain numbers, not
# formulas, so Excel's row insert does not renumber the rows that follow.
# Recompute them for every row from the insertion point through the end
# of the sheet (row 179 after the insert).
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 164; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = ($r - 2)
    $ws.Cells.Item($r, 4).Value = ($r - 1)
}

# --- Occurrence count bump for the "人" row, picked up by the re-crawl ---
$ws.Cells.Item(6, 7).Value = 39
